$p = $ppt.ActivePresentation

# Append a new slide at the end using the "Title and Content" layout
# (ppLayoutText = 2), matching the layout used elsewhere in the deck
# (ppt/slideLayouts/slideLayout2.xml - "Titel og indholdsobjekt").
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# Title placeholder: "Evt. noget 2. ordens højpasfilter?"
# Built out of three runs (to mirror the authored run/spell-check
# boundaries) using InsertAfter, which appends as a distinct run.
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Evt. noget 2. ordens "
$title.LanguageID = "da-DK"
[void]$title.InsertAfter("højpasfilter")
[void]$title.InsertAfter("?")

# Keep the body/content placeholder empty but give it the Danish
# language context, same as the rest of the deck.
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.LanguageID = "da-DK"
